$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 10:22"

# Row 24 - Austria
$ws.Range("B24").Value = 14940
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 11694
$ws.Range("E24").Value = 2736
$ws.Range("F24").Value = 169

# Row 34 - Polonia
$ws.Range("B34").Value = 10346
$ws.Range("C34").Value = 177
$ws.Range("E34").Value = 8171
$ws.Range("G34").Value = 9
$ws.Range("H34").Value = 435

# Row 37 - Dinamarca
$ws.Range("B37").Value = 8073
$ws.Range("C37").Value = 161
$ws.Range("E37").Value = 2602

# Row 45 - Filipinas
$ws.Range("B45").Value = 6981
$ws.Range("C45").Value = 271
$ws.Range("D45").Value = 722
$ws.Range("E45").Value = 5797
$ws.Range("G45").Value = 16
$ws.Range("H45").Value = 462

# Row 60 - Moldavia
$ws.Range("D60").Value = 661
$ws.Range("E60").Value = 2041
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 76

# Row 64 - Kazajistan
$ws.Range("B64").Value = 2207
$ws.Range("C64").Value = 72
$ws.Range("D64").Value = 536
$ws.Range("E64").Value = 1651
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 20

$wb.Save()
